$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 0.5232980847358704
$ws.Cells.Item(2, 3).Value = 2.709017276763916
$ws.Cells.Item(2, 4).Value = 0.4522778391838074
$ws.Cells.Item(2, 5).Value = 3.056704521179199
$ws.Cells.Item(2, 6).Value = 0.4223107695579529
$ws.Cells.Item(2, 7).Value = 3.020390510559082
$ws.Cells.Item(2, 8).Value = 0.3550580143928528
$ws.Cells.Item(2, 9).Value = 3.237790584564209
$ws.Cells.Item(2, 10).Value = 0.54347825050354
$ws.Cells.Item(2, 11).Value = 2.821513414382935
$ws.Cells.Item(2, 12).Value = 0.4821583330631256
$ws.Cells.Item(2, 13).Value = 3.039693117141724
$ws.Cells.Item(2, 14).Value = 0.5365927815437317
$ws.Cells.Item(2, 15).Value = 2.578890562057495
$ws.Cells.Item(2, 16).Value = 0.5165425539016724
$ws.Cells.Item(2, 17).Value = 2.82404613494873
$ws.Cells.Item(2, 18).Value = 0.5649142265319824
$ws.Cells.Item(2, 19).Value = 2.681746244430542
$ws.Cells.Item(2, 20).Value = 0.5009527206420898
$ws.Cells.Item(2, 21).Value = 2.894340753555298
$ws.Cells.Item(2, 22).Value = 0.5574657917022705
$ws.Cells.Item(2, 23).Value = 2.610142230987549
$ws.Cells.Item(2, 24).Value = 0.5297072529792786
$ws.Cells.Item(2, 25).Value = 2.833125591278076
$ws.Cells.Item(2, 26).Value = 0.5851377248764038
$ws.Cells.Item(2, 27).Value = 2.72585916519165
$ws.Cells.Item(2, 28).Value = 0.5334748029708862
$ws.Cells.Item(2, 29).Value = 2.744796514511108
$ws.Cells.Item(2, 30).Value = 0.5152000784873962
$ws.Cells.Item(2, 31).Value = 2.862178325653076
$ws.Cells.Item(2, 32).Value = 42.6402473449707

# Row 3
$ws.Cells.Item(3, 2).Value = 0.824701189994812
$ws.Cells.Item(3, 3).Value = 0.5300958752632141
$ws.Cells.Item(3, 4).Value = 0.6448986530303955
$ws.Cells.Item(3, 5).Value = 1.110059261322021
$ws.Cells.Item(3, 6).Value = 0.7622120380401611
$ws.Cells.Item(3, 7).Value = 0.7923889756202698
$ws.Cells.Item(3, 8).Value = 0.61168372631073
$ws.Cells.Item(3, 9).Value = 1.142401099205017
$ws.Cells.Item(3, 10).Value = 0.8251775503158569
$ws.Cells.Item(3, 11).Value = 0.5393660664558411
$ws.Cells.Item(3, 12).Value = 0.8227957487106323
$ws.Cells.Item(3, 13).Value = 0.5907964706420898
$ws.Cells.Item(3, 14).Value = 0.8232721090316772
$ws.Cells.Item(3, 15).Value = 0.6082631349563599
$ws.Cells.Item(3, 16).Value = 0.822839081287384
$ws.Cells.Item(3, 17).Value = 0.539411723613739
$ws.Cells.Item(3, 18).Value = 0.7875887751579285
$ws.Cells.Item(3, 19).Value = 0.6492300033569336
$ws.Cells.Item(3, 20).Value = 0.8213667273521423
$ws.Cells.Item(3, 21).Value = 0.5737504959106445
$ws.Cells.Item(3, 22).Value = 0.7985449433326721
$ws.Cells.Item(3, 23).Value = 0.6190176606178284
$ws.Cells.Item(3, 25).Value = 0.5296198129653931
$ws.Cells.Item(3, 26).Value = 0.8292482495307922
$ws.Cells.Item(3, 27).Value = 0.5291793942451477
$ws.Cells.Item(3, 28).Value = 0.8248744010925293
$ws.Cells.Item(3, 29).Value = 0.6159378290176392
$ws.Cells.Item(3, 30).Value = 0.8282955288887024
$ws.Cells.Item(3, 31).Value = 0.5404333472251892
$ws.Cells.Item(3, 32).Value = 9.909952163696289

# Row 4
$ws.Cells.Item(4, 2).Value = 0.8359605073928833
$ws.Cells.Item(4, 3).Value = 0.3564878106117249
$ws.Cells.Item(4, 4).Value = 0.6538628339767456
$ws.Cells.Item(4, 5).Value = 0.9263015985488892
$ws.Cells.Item(4, 6).Value = 0.7723886966705322
$ws.Cells.Item(4, 7).Value = 0.5942257046699524
$ws.Cells.Item(4, 8).Value = 0.6293521523475647
$ws.Cells.Item(4, 9).Value = 0.9403344392776489
$ws.Cells.Item(4, 10).Value = 0.8347912430763245
$ws.Cells.Item(4, 11).Value = 0.3576248586177826
$ws.Cells.Item(4, 12).Value = 0.8379958271980286
$ws.Cells.Item(4, 13).Value = 0.3619590997695923
$ws.Cells.Item(4, 14).Value = 0.8434956073760986
$ws.Cells.Item(4, 15).Value = 0.4229556322097778
$ws.Cells.Item(4, 16).Value = 0.8369132280349731
$ws.Cells.Item(4, 17).Value = 0.3551364243030548
$ws.Cells.Item(4, 18).Value = 0.7977654337882996
$ws.Cells.Item(4, 19).Value = 0.4765941202640533
$ws.Cells.Item(4, 20).Value = 0.833968460559845
$ws.Cells.Item(4, 21).Value = 0.3619384467601776
$ws.Cells.Item(4, 22).Value = 0.8115364909172058
$ws.Cells.Item(4, 23).Value = 0.450566440820694
$ws.Cells.Item(4, 24).Value = 0.8365668058395386
$ws.Cells.Item(4, 25).Value = 0.3543886244297028
$ws.Cells.Item(4, 26).Value = 0.8367832899093628
$ws.Cells.Item(4, 27).Value = 0.3551085591316223
$ws.Cells.Item(4, 28).Value = 0.8426294922828674
$ws.Cells.Item(4, 29).Value = 0.4256309866905212
$ws.Cells.Item(4, 30).Value = 0.8381257653236389
$ws.Cells.Item(4, 31).Value = 0.3539482653141022
$ws.Cells.Item(4, 32).Value = 7.093203067779541

# Row 5
$ws.Cells.Item(5, 2).Value = 0.8409838676452637
$ws.Cells.Item(5, 3).Value = 0.3334563970565796
$ws.Cells.Item(5, 4).Value = 0.6613546013832092
$ws.Cells.Item(5, 5).Value = 0.9001128077507019
$ws.Cells.Item(5, 6).Value = 0.7780616879463196
$ws.Cells.Item(5, 7).Value = 0.5657532215118408
$ws.Cells.Item(5, 8).Value = 0.6340290904045105
$ws.Cells.Item(5, 9).Value = 0.9101560711860657
$ws.Cells.Item(5, 10).Value = 0.8415468335151672
$ws.Cells.Item(5, 11).Value = 0.334507554769516
$ws.Cells.Item(5, 12).Value = 0.8427160978317261
$ws.Cells.Item(5, 13).Value = 0.337427169084549
$ws.Cells.Item(5, 14).Value = 0.8464836478233337
$ws.Cells.Item(5, 15).Value = 0.4020196497440338
$ws.Cells.Item(5, 16).Value = 0.8433223366737366
$ws.Cells.Item(5, 17).Value = 0.3322597146034241
$ws.Cells.Item(5, 18).Value = 0.8012298345565796
$ws.Cells.Item(5, 19).Value = 0.452781468629837
$ws.Cells.Item(5, 20).Value = 0.8417201042175293
$ws.Cells.Item(5, 21).Value = 0.3363668024539948
$ws.Cells.Item(5, 22).Value = 0.814048171043396
$ws.Cells.Item(5, 23).Value = 0.4252075552940369
$ws.Cells.Item(5, 24).Value = 0.8422830700874329
$ws.Cells.Item(5, 25).Value = 0.3322798609733582
$ws.Cells.Item(5, 26).Value = 0.8425861597061157
$ws.Cells.Item(5, 27).Value = 0.3339123427867889
$ws.Cells.Item(5, 28).Value = 0.8451845049858093
$ws.Cells.Item(5, 29).Value = 0.4042204916477203
$ws.Cells.Item(5, 30).Value = 0.8415468335151672
$ws.Cells.Item(5, 31).Value = 0.3322315812110901
$ws.Cells.Item(5, 32).Value = 6.732691287994385

# Row 6
$ws.Cells.Item(6, 2).Value = 0.8430192470550537
$ws.Cells.Item(6, 3).Value = 0.3289691209793091
$ws.Cells.Item(6, 4).Value = 0.6643859148025513
$ws.Cells.Item(6, 5).Value = 0.8918812870979309
$ws.Cells.Item(6, 6).Value = 0.7792742252349854
$ws.Cells.Item(6, 7).Value = 0.5582778453826904
$ws.Cells.Item(6, 8).Value = 0.63831627368927
$ws.Cells.Item(6, 9).Value = 0.9025792479515076
$ws.Cells.Item(6, 10).Value = 0.8424995541572571
$ws.Cells.Item(6, 11).Value = 0.3293662071228027
$ws.Cells.Item(6, 12).Value = 0.84340900182724
$ws.Cells.Item(6, 13).Value = 0.3339523375034332
$ws.Cells.Item(6, 14).Value = 0.8479126691818237
$ws.Cells.Item(6, 15).Value = 0.3986920416355133
$ws.Cells.Item(6, 16).Value = 0.8432790637016296
$ws.Cells.Item(6, 17).Value = 0.3287070095539093
$ws.Cells.Item(6, 18).Value = 0.803178608417511
$ws.Cells.Item(6, 19).Value = 0.4481481909751892
$ws.Cells.Item(6, 20).Value = 0.8416767716407776
$ws.Cells.Item(6, 21).Value = 0.3331458270549774
$ws.Cells.Item(6, 22).Value = 0.816126823425293
$ws.Cells.Item(6, 23).Value = 0.4213265776634216
$ws.Cells.Item(6, 24).Value = 0.8436254858970642
$ws.Cells.Item(6, 25).Value = 0.3274954557418823
$ws.Cells.Item(6, 26).Value = 0.8424995541572571
$ws.Cells.Item(6, 27).Value = 0.3298241794109344
$ws.Cells.Item(6, 28).Value = 0.8474796414375305
$ws.Cells.Item(6, 29).Value = 0.3992651700973511
$ws.Cells.Item(6, 30).Value = 0.8381257653236389
$ws.Cells.Item(6, 31).Value = 0.3282879292964935
